{"js": "// Update the three-digit \u00f7 one-digit answer table with newly generated values.\n// Each cell's text is located by its exact current content and replaced with\n// the new computed text.\n\nconst replacements = [\n  [\"947\u00f79=105, 2\", \"325\u00f76=54, 1\"],\n  [\"245\u00f77=35, 0\", \"740\u00f78=92, 4\"],\n  [\"847\u00f73=282, 1\", \"498\u00f76=83, 0\"],\n  [\"754\u00f76=125, 4\", \"889\u00f72=444, 1\"],\n  [\"489\u00f74=122, 1\", \"555\u00f72=277, 1\"],\n  [\"969\u00f77=138, 3\", \"297\u00f78=37, 1\"],\n  [\"897\u00f76=149, 3\", \"728\u00f76=121, 2\"],\n  [\"685\u00f76=114, 1\", \"672\u00f75=134, 2\"],\n  [\"747\u00f76=124, 3\", \"278\u00f75=55, 3\"],\n  [\"816\u00f73=272, 0\", \"867\u00f76=144, 3\"],\n  [\"983\u00f75=196, 3\", \"432\u00f78=54, 0\"],\n  [\"680\u00f78=85, 0\", \"800\u00f72=400, 0\"],\n  [\"637\u00f72=318, 1\", \"482\u00f76=80, 2\"],\n  [\"914\u00f78=114, 2\", \"557\u00f74=139, 1\"],\n  [\"370\u00f78=46, 2\", \"262\u00f77=37, 3\"],\n  [\"739\u00f72=369, 1\", \"622\u00f77=88, 6\"],\n  [\"852\u00f79=94, 6\", \"457\u00f74=114, 1\"],\n  [\"852\u00f77=121, 5\", \"849\u00f75=169, 4\"],\n  [\"246\u00f72=123, 0\", \"319\u00f78=39, 7\"],\n  [\"693\u00f74=173, 1\", \"958\u00f77=136, 6\"],\n  [\"431\u00f77=61, 4\", \"526\u00f75=105, 1\"],\n  [\"973\u00f79=108, 1\", \"252\u00f77=36, 0\"],\n  [\"929\u00f76=154, 5\", \"683\u00f79=75, 8\"],\n  [\"425\u00f79=47, 2\", \"598\u00f76=99, 4\"],\n  [\"975\u00f72=487, 1\", \"421\u00f78=52, 5\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit \u00f7 one-digit answer table with newly generated values.\n# Each cell's text is replaced in place by matching its exact current content\n# and substituting the new computed text (so ordering/formatting of the run\n# is preserved).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"947\u00f79=105, 2\"; New = \"325\u00f76=54, 1\"},\n    @{Old = \"245\u00f77=35, 0\";  New = \"740\u00f78=92, 4\"},\n    @{Old = \"847\u00f73=282, 1\"; New = \"498\u00f76=83, 0\"},\n    @{Old = \"754\u00f76=125, 4\"; New = \"889\u00f72=444, 1\"},\n    @{Old = \"489\u00f74=122, 1\"; New = \"555\u00f72=277, 1\"},\n    @{Old = \"969\u00f77=138, 3\"; New = \"297\u00f78=37, 1\"},\n    @{Old = \"897\u00f76=149, 3\"; New = \"728\u00f76=121, 2\"},\n    @{Old = \"685\u00f76=114, 1\"; New = \"672\u00f75=134, 2\"},\n    @{Old = \"747\u00f76=124, 3\"; New = \"278\u00f75=55, 3\"},\n    @{Old = \"816\u00f73=272, 0\"; New = \"867\u00f76=144, 3\"},\n    @{Old = \"983\u00f75=196, 3\"; New = \"432\u00f78=54, 0\"},\n    @{Old = \"680\u00f78=85, 0\";  New = \"800\u00f72=400, 0\"},\n    @{Old = \"637\u00f72=318, 1\"; New = \"482\u00f76=80, 2\"},\n    @{Old = \"914\u00f78=114, 2\"; New = \"557\u00f74=139, 1\"},\n    @{Old = \"370\u00f78=46, 2\";  New = \"262\u00f77=37, 3\"},\n    @{Old = \"739\u00f72=369, 1\"; New = \"622\u00f77=88, 6\"},\n    @{Old = \"852\u00f79=94, 6\";  New = \"457\u00f74=114, 1\"},\n    @{Old = \"852\u00f77=121, 5\"; New = \"849\u00f75=169, 4\"},\n    @{Old = \"246\u00f72=123, 0\"; New = \"319\u00f78=39, 7\"},\n    @{Old = \"693\u00f74=173, 1\"; New = \"958\u00f77=136, 6\"},\n    @{Old = \"431\u00f77=61, 4\";  New = \"526\u00f75=105, 1\"},\n    @{Old = \"973\u00f79=108, 1\"; New = \"252\u00f77=36, 0\"},\n    @{Old = \"929\u00f76=154, 5\"; New = \"683\u00f79=75, 8\"},\n    @{Old = \"425\u00f79=47, 2\";  New = \"598\u00f76=99, 4\"},\n    @{Old = \"975\u00f72=487, 1\"; New = \"421\u00f78=52, 5\"}\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 0\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
